$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header labels (D1, E1) to lowercase "ing" ---
$ws.Cells.Item(1,4).Value2 = "metadata4ing_IRI"
$ws.Cells.Item(1,5).Value2 = "metadata4ing_DESC"

# --- Add new column F: header + data ---
# Copy formatting from E1 (header style) to F1
$ws.Cells.Item(1,5).Copy()
$ws.Cells.Item(1,6).PasteSpecial(-4122)
$ws.Cells.Item(1,6).Value2 = "metadata4ing_DEF"

# Copy formatting from E2 (data style) to F2
$ws.Cells.Item(2,5).Copy()
$ws.Cells.Item(2,6).PasteSpecial(-4122)
$ws.Cells.Item(2,6).Value2 = "[locstr('Object that helps an agent perform an action (e.g., John wrote a book with a pen), i.e., an object of the relation schema:instrument or m4i:hasTool', 'en')]"

$excel.CutCopyMode = 0
